$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "EventBridge " + "rule"  ->  single run "EventBridge rule" ---
# (the paragraph's endParaRPr must be preserved, so force a genuine
#  text change - via a throwaway value - to make the two runs coalesce
#  back into one run when we set the final text.)
$tr1 = $s.Shapes.Item("Google Shape;68;p13").TextFrame.TextRange
$tr1.Text = "zzz__placeholder__zzz"
$tr1.Text = "EventBridge rule"

# --- 2) "System Manager Automation runbook" (1 run) ---
#        -> "Systems " (new run) + "Manager Automation runbook" ---
$tr2 = $s.Shapes.Item("Google Shape;71;p13").TextFrame.TextRange
$c2 = $tr2.Characters(1, 7)
$c2.Text = "Systems "

# --- 3) "Systems Manager Command " + "d" + "ocument" ---
#        -> single run "Systems Manager Command document" ---
$tr3 = $s.Shapes.Item("Google Shape;72;p13").TextFrame.TextRange
$tr3.Text = "zzz__placeholder__zzz"
$tr3.Text = "Systems Manager Command document"

# --- 4) "Automation " + "execution " + "r" + "ole" ---
#        -> single run "Automation execution IAM role" ---
#        (leave the preceding "Systems Manager " run untouched) ---
$tr4 = $s.Shapes.Item("Google Shape;75;p13").TextFrame.TextRange
$c4 = $tr4.Characters(17, 26)
$c4.Text = "Automation execution IAM role"

# --- 5) AWS IAM text box (inside Group 9) - drop the stray endParaRPr ---
$g5 = $s.Shapes.Item("Group 9")
$tr5 = $g5.GroupItems.Item("TextBox 12").TextFrame.TextRange
[void]$tr5.Delete()
[void]$tr5.InsertAfter("AWS IAM")

# --- 6) "New Relic license key" text box - drop the stray endParaRPr ---
$tr6 = $s.Shapes.Item("TextBox 6").TextFrame.TextRange
[void]$tr6.Delete()
[void]$tr6.InsertAfter("New Relic license key")
